# SBERDOMA-1030 - refactor with global mappers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Complete" -> "Completed" (header cell C1 uses this shared string)
$ws.Range("C1").Value = "Completed"

# Column widths (characters) tuned so the stored xlsx <col> width matches
# the target widths as closely as the engine's rounding allows:
#   B: 21.36 -> 23.98
#   C: 11.63 -> 16.96 (now a custom width)
#   D: 14.77 -> 16.22
#   E: 11.63 -> 14.62 (now a custom width)
$ws.Columns.Item(2).ColumnWidth = 23.166666666666668
$ws.Columns.Item(3).ColumnWidth = 16.166666666666668
$ws.Columns.Item(4).ColumnWidth = 15.333333333333334
$ws.Columns.Item(5).ColumnWidth = 13.833333333333334

# Row heights for data rows 2 and 3: 36.9 -> 25.1
$ws.Rows.Item(2).RowHeight = 25.1
$ws.Rows.Item(3).RowHeight = 25.1

# Move the active selection from D9 to E9
[void]$ws.Range("E9").Select()
